$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-26 Thursday", 2)

$d.Content.Find.Execute("21×50=", $true, $false, $false, $false, $false, $true, 1, $false, "78×16=", 2)
$d.Content.Find.Execute("44×18=", $true, $false, $false, $false, $false, $true, 1, $false, "31×58=", 2)
$d.Content.Find.Execute("64×15=", $true, $false, $false, $false, $false, $true, 1, $false, "13×53=", 2)
$d.Content.Find.Execute("61×94=", $true, $false, $false, $false, $false, $true, 1, $false, "60×59=", 2)
$d.Content.Find.Execute("44×33=", $true, $false, $false, $false, $false, $true, 1, $false, "81×75=", 2)

$d.Content.Find.Execute("12×50=", $true, $false, $false, $false, $false, $true, 1, $false, "24×13=", 2)
$d.Content.Find.Execute("61×43=", $true, $false, $false, $false, $false, $true, 1, $false, "60×49=", 2)
$d.Content.Find.Execute("45×38=", $true, $false, $false, $false, $false, $true, 1, $false, "73×91=", 2)
$d.Content.Find.Execute("75×88=", $true, $false, $false, $false, $false, $true, 1, $false, "22×51=", 2)
$d.Content.Find.Execute("67×22=", $true, $false, $false, $false, $false, $true, 1, $false, "82×55=", 2)

$d.Content.Find.Execute("59×84=", $true, $false, $false, $false, $false, $true, 1, $false, "75×39=", 2)
$d.Content.Find.Execute("14×89=", $true, $false, $false, $false, $false, $true, 1, $false, "36×64=", 2)
$d.Content.Find.Execute("59×61=", $true, $false, $false, $false, $false, $true, 1, $false, "15×74=", 2)
$d.Content.Find.Execute("72×82=", $true, $false, $false, $false, $false, $true, 1, $false, "56×89=", 2)
$d.Content.Find.Execute("35×78=", $true, $false, $false, $false, $false, $true, 1, $false, "85×28=", 2)

$d.Content.Find.Execute("86×31=", $true, $false, $false, $false, $false, $true, 1, $false, "74×18=", 2)
$d.Content.Find.Execute("41×19=", $true, $false, $false, $false, $false, $true, 1, $false, "97×20=", 2)
$d.Content.Find.Execute("79×32=", $true, $false, $false, $false, $false, $true, 1, $false, "13×70=", 2)
$d.Content.Find.Execute("32×81=", $true, $false, $false, $false, $false, $true, 1, $false, "71×68=", 2)
$d.Content.Find.Execute("29×60=", $true, $false, $false, $false, $false, $true, 1, $false, "43×83=", 2)

$d.Content.Find.Execute("20×43=", $true, $false, $false, $false, $false, $true, 1, $false, "88×65=", 2)
$d.Content.Find.Execute("58×44=", $true, $false, $false, $false, $false, $true, 1, $false, "52×12=", 2)
$d.Content.Find.Execute("68×14=", $true, $false, $false, $false, $false, $true, 1, $false, "30×57=", 2)
$d.Content.Find.Execute("14×36=", $true, $false, $false, $false, $false, $true, 1, $false, "84×96=", 2)
$d.Content.Find.Execute("22×18=", $true, $false, $false, $false, $false, $true, 1, $false, "96×65=", 2)
